$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume figures per the Sat May 13 2023 GitHub Actions refresh.
# Force column D (Price) and E (Volume(1h)) to remain plain text so values such as
# "318.49", "4.470", or "0.07170" are not silently reinterpreted as numbers and lose
# their original formatting (trailing zeros, dot-grouped thousands, etc.).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.419.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +2.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.49"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4361"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3723"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07347"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8725"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.967.19"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +10.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.478"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.681"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07170"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.24"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.031"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008993"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.23%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.446.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.250"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.15"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.168.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.898"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.55"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.269"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.919"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.198"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7582"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.470"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.151"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01956"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05243"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5161"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.795"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.36%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.524"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.460"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.66"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4632"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.879"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.25%  "
